# Apply updated vm_pu values for the "case with 380 kV" run.
# Columns B..F and I..M of rows 2-25 receive new voltage-magnitude values
# (column A index, G and N, and the header row 1 are unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row: @(RowNumber, B, C, D, E, F, I, J, K, L, M)
$data = @(
    @(2, 1.02, 1.057742643088169, 1.055265475451534, 1.063305949821053, 1.072666294045031, 1.048497942368824, 1.062737480654799, 1.058006408157197, 1.066024946841664, 1.075360202523922),
    @(3, 1.02, 1.05888689826598, 1.056120190283748, 1.064334219556289, 1.073797905950379, 1.048832845641869, 1.063533389887483, 1.058674424707046, 1.06686766622385, 1.076307826669176),
    @(4, 1.02, 1.059627357309382, 1.056673222483377, 1.064999939573589, 1.074530671708417, 1.049048396441891, 1.064047873578672, 1.059106006876108, 1.067412712781813, 1.076920948638579),
    @(5, 1.02, 1.059938658940513, 1.056905711024509, 1.065279894779086, 1.074838855317227, 1.049138737794471, 1.064264037615963, 1.059287283906702, 1.067641790886076, 1.077178692379642),
    @(6, 1.02, 1.059990928619758, 1.056944746507281, 1.065326905567316, 1.074890608267277, 1.049153890310421, 1.064300325191844, 1.059317711731517, 1.067680250626315, 1.077221967948303),
    @(7, 1.02, 1.059631516886276, 1.056676329031114, 1.065003680008628, 1.074534789166095, 1.049049604672666, 1.064050762463592, 1.059108429737676, 1.067415773968073, 1.076924392672788),
    @(8, 1.02, 1.058129339738296, 1.055554335539, 1.06365338378799, 1.073048616527131, 1.04861136348502, 1.063006570308762, 1.058232306099749, 1.066309799527064, 1.075680467914101),
    @(9, 1.02, 1.05548265348719, 1.053577048780785, 1.061276745052431, 1.070433891671437, 1.047830284881648, 1.061162554552437, 1.056683334283536, 1.064359013159382, 1.073488086938376),
    @(10, 1.02, 1.0537183718387, 1.052258732058412, 1.059694159173466, 1.068693478064747, 1.047303615240334, 1.05993049195408, 1.055647228444097, 1.06305718383358, 1.072026194208775),
    @(11, 1.02, 1.052954447290713, 1.051687856299236, 1.059009314651614, 1.067940503135829, 1.04707414830438, 1.059396345447071, 1.05519776031168, 1.06249316358021, 1.071393101094776),
    @(12, 1.02, 1.052670693653543, 1.051475802099835, 1.058754996417077, 1.067660909628143, 1.046988701185288, 1.059197840681059, 1.055030683194492, 1.062283612872215, 1.071157928985121),
    @(13, 1.02, 1.052731559702917, 1.05152128869749, 1.058809545674457, 1.067720879099393, 1.047007039519243, 1.059240425114225, 1.055066527448699, 1.06232856441133, 1.071208374804029),
    @(14, 1.02, 1.052930992093671, 1.051670327946396, 1.058988291333834, 1.067917389916604, 1.047067089562846, 1.059379939003177, 1.05518395220954, 1.062475843044817, 1.071373661965922),
    @(15, 1.02, 1.053053869194497, 1.051762155210719, 1.059098430837269, 1.068038479275202, 1.047104060161173, 1.059465884984527, 1.055256284954519, 1.062566579834789, 1.071475499084662),
    @(16, 1.02, 1.053769071335272, 1.05229661843476, 1.059739618979323, 1.06874346387275, 1.047318814363389, 1.059965927648074, 1.055677040702427, 1.063094609227936, 1.072068208700336),
    @(17, 1.02, 1.054217703524464, 1.0526318633887, 1.060141933005299, 1.06918585184855, 1.047453144875125, 1.060279415395715, 1.055940747902787, 1.063425742233347, 1.072439976941479),
    @(18, 1.02, 1.054479385201978, 1.052827402641444, 1.060376637063758, 1.069443950640682, 1.047531361028708, 1.060462204305982, 1.056094484034281, 1.063618855710694, 1.072656815005735),
    @(19, 1.02, 1.054568612286542, 1.052894075890942, 1.060456672065651, 1.06953196608643, 1.047558007559305, 1.060524519898648, 1.056146890553712, 1.063684697205409, 1.072730749830938),
    @(20, 1.02, 1.054169569326476, 1.052595895136248, 1.060098764270213, 1.069138381461746, 1.047438746609186, 1.060245787666807, 1.055912462873945, 1.063390217984917, 1.072400090565063),
    @(21, 1.02, 1.052872264181078, 1.051626439755826, 1.058935653442426, 1.067859519763686, 1.047049412200045, 1.05933885837644, 1.055149376987478, 1.062432474499874, 1.071324989383373),
    @(22, 1.02, 1.052056607653305, 1.051016872655342, 1.058204726759069, 1.067055997467426, 1.04680339042755, 1.058768063064962, 1.05466887331453, 1.061830022480809, 1.070648954321431),
    @(23, 1.02, 1.052489001999541, 1.051340018787316, 1.058592170155003, 1.067481908016284, 1.046933928054761, 1.059070706935461, 1.05492366588058, 1.062149420453979, 1.071007340706315),
    @(24, 1.02, 1.054191319091735, 1.052612147651829, 1.060118270236234, 1.069159831096326, 1.047445252987638, 1.060260982785035, 1.055925243907414, 1.063406269958622, 1.072418113526351),
    @(25, 1.02, 1.05616685038301, 1.054088247130167, 1.061890837888071, 1.071109376884251, 1.048033260820294, 1.061639754651754, 1.057084389469652, 1.064863567092965, 1.074054922443649)
)

foreach ($rowData in $data) {
    $r = $rowData[0]
    $ws.Range("B$r").Value = $rowData[1]
    $ws.Range("C$r").Value = $rowData[2]
    $ws.Range("D$r").Value = $rowData[3]
    $ws.Range("E$r").Value = $rowData[4]
    $ws.Range("F$r").Value = $rowData[5]
    $ws.Range("I$r").Value = $rowData[6]
    $ws.Range("J$r").Value = $rowData[7]
    $ws.Range("K$r").Value = $rowData[8]
    $ws.Range("L$r").Value = $rowData[9]
    $ws.Range("M$r").Value = $rowData[10]
}
